$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: A4 gets a new shared string "Get Bucket List".
$ws.Range("A4").Value = "Get Bucket List"

# The authored workbook registers a custom table style ("MySqlDefault",
# built from a bold/shaded "whole table" dxf and a plain "header row" dxf)
# in styles.xml even though no live Excel Table/ListObject remains on the
# sheet. Driving that through a transient ListObject and then unlisting it
# is the most faithful way to reproduce that artifact via the object model.
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:A4"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "MySqlDefault"
$lo.TableStyle = "TableStyleMedium2"
$lo.Unlist()

# Match the final selection left behind by the edit (B4).
$ws.Range("B4").Select()
